$d = $word.ActiveDocument

# NOTE: the replacement text also occurs verbatim as a substring of the
# "Version history and licensing" paragraph further down the document
# (e.g. "...created 08/23 by Zoë Gemmell, Isabella Lewis, Akshat
# Srivastava as part of..."), so the Find/Replace below is scoped to each
# specific paragraph's own Range rather than run against $d.Content -
# otherwise it would also rewrite (and collapse the whitespace-handling
# of) that unrelated, already-single-run paragraph.

# --- 1. Collapse the Title paragraph's word-by-word runs into one run ---
$d.Paragraphs(1).Range.Find.Execute(
    "Questions: Solving exponential equations", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Questions: Solving exponential equations", 2) | Out-Null

# --- 2. Collapse the Author paragraph's word-by-word runs into one run ---
$d.Paragraphs(2).Range.Find.Execute(
    "Zoë Gemmell, Isabella Lewis, Akshat Srivastava", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Zoë Gemmell, Isabella Lewis, Akshat Srivastava", 2) | Out-Null

# --- 3. Collapse the Abstract paragraph's word-by-word runs into one run ---
$d.Paragraphs(4).Range.Find.Execute(
    "A selection of questions for the study guide on solving equations involving indices.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A selection of questions for the study guide on solving equations involving indices.", 2) | Out-Null

# --- 4. Re-order the <m:dPr> children of every math delimiter so that
#        <m:sepChr/> sits between <m:begChr/> and <m:endChr/> (matches the
#        canonical CT_OMathDPr sequence: begChr, sepChr, endChr, grow). The
#        Word object model has no direct property for this low-level child
#        ordering, so each affected equation's OOXML is round-tripped
#        through WordOpenXML / InsertXML with the delimiter properties
#        rewritten in place. ---
$dPrPattern = '<m:dPr><m:begChr m:val="([^"]*)"\s*/><m:endChr m:val="([^"]*)"\s*/><m:sepChr m:val=""\s*/><m:grow\s*/></m:dPr>'
$dPrReplacement = '<m:dPr><m:begChr m:val="$1"/><m:sepChr m:val=""/><m:endChr m:val="$2"/><m:grow/></m:dPr>'

for ($i = 1; $i -le $d.OMaths.Count; $i++) {
    $om = $d.OMaths.Item($i)
    $r = $om.Range
    $xml = $r.WordOpenXML
    if ([regex]::IsMatch($xml, $dPrPattern)) {
        $newXml = [regex]::Replace($xml, $dPrPattern, $dPrReplacement)
        $r.InsertXML($newXml)
    }
}
